$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name (date rolled forward from 10-06 to 10-07)
$ws.Name = "Through 2022-10-07"

# Update header label for the "Total"/current-year column
$ws.Range("I1").Value = "2022 (through 10-07)"

# Update data values for new day's data (2022-10-15 data add)
$ws.Range("I11").Value = 22
$ws.Range("H12").Value = 202
$ws.Range("H14").Value = 1849
$ws.Range("I14").Value = 1303
